$d = $word.ActiveDocument

$newText = "Osallistut maailmanlaajuiseen kampanjaan tarkkaillaksesi ja tallentaaksesi himmeimpiä näkyvissä olevia tähtiä keinona mitata valonsaastetta tietyssä paikassa. Paikallistamalla ja tarkkailemalla Herkuleen tähtikuvio miten valosaaste syntyy kunkin taajaman tai muun ihmisen toiminnan valoista. Antamasi tiedot päivittyvät heti verkossa olevaan tietokantaan, ja näin saadaan käsitys siitä minkä verran taivaan tähdistä on missäkin nähtävissä."

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Osallistut maailmanlaajuiseen*") {
        $r = $p.Range
        # Exclude the trailing paragraph mark so the paragraph itself is kept.
        $r.End = $r.End - 1
        # Remove all existing runs/formatting in the paragraph...
        $r.Text = ""
        # ...then append fresh, unformatted text as a brand-new run.
        $r.InsertAfter($newText)
        break
    }
}
